$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): values are plain text (some contain multiple "." as thousands
# separators), so force text format before assigning to avoid Excel coercing them
# into numbers, then restore the default "Normal" style so no extra cell formatting
# is introduced.
$priceUpdates = @{
    "D2" = "67.795.74"
    "D3" = "2.491.69"
    "D5" = "586.84"
    "D6" = "177.15"
    "D13" = "2.947.46"
    "D15" = "67.719.08"
    "D17" = "2.500.56"
    "D19" = "10.96"
    "D20" = "350.31"
    "D21" = "4.11"
    "D23" = "70.84"
    "D24" = "4.26"
    "D26" = "9.10"
    "D28" = "0.999"
    "D30" = "504.37"
    "D36" = "162.94"
    "D45" = "144.67"
    "D49" = "0.0741"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Volume(1h) column (E): percentage strings padded with two spaces on each side.
$volumeUpdates = @{
    "E2" = "  +1.07%  "
    "E3" = "  +0.37%  "
    "E4" = "  +0.07%  "
    "E5" = "  +0.29%  "
    "E6" = "  +3.40%  "
    "E7" = "  -0.04%  "
    "E8" = "  +0.33%  "
    "E9" = "  +3.86%  "
    "E10" = "  +0.14%  "
    "E12" = "  +0.13%  "
    "E13" = "  +0.45%  "
    "E14" = "  +1.00%  "
    "E16" = "  +0.87%  "
    "E17" = "  +1.02%  "
    "E18" = "  +1.64%  "
    "E19" = "  -0.07%  "
    "E20" = "  +0.11%  "
    "E21" = "  +2.31%  "
    "E22" = "  -0.24%  "
    "E23" = "  +3.44%  "
    "E24" = "  +0.85%  "
    "E25" = "  -2.60%  "
    "E26" = "  -1.54%  "
    "E27" = "  +0.31%  "
    "E28" = "  -0.03%  "
    "E29" = "  +0.36%  "
    "E30" = "  -1.31%  "
    "E31" = "  +1.65%  "
    "E32" = "  +2.52%  "
    "E33" = "  +0.60%  "
    "E34" = "  +0.01%  "
    "E35" = "  +3.73%  "
    "E36" = "  +1.93%  "
    "E38" = "  +0.50%  "
    "E39" = "  +0.44%  "
    "E41" = "  +3.68%  "
    "E42" = "  +0.27%  "
    "E43" = "  +0.88%  "
    "E44" = "  +1.82%  "
    "E45" = "  +1.30%  "
    "E46" = "  +2.26%  "
    "E47" = "  +0.00%  "
    "E48" = "  +1.59%  "
    "E49" = "  +1.64%  "
    "E50" = "  +1.21%  "
    "E51" = "  +0.53%  "
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

Write-Host "Updated $($priceUpdates.Count) price cells and $($volumeUpdates.Count) volume cells"
